# Generate Report for Archive
$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Update status text from "Ready for handoff" to "In Translation"
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# Narrow the status columns to fit the new (shorter) text.
# Target stored OOXML column width is 13.4101845877511 characters; the
# COM ColumnWidth setter here quantizes the stored width to the nearest
# 1/6th of a character (stored = round(ColumnWidth*6)/6 + 5/6), so 12.5
# is the input value that lands closest to the target (13.3333...).
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
